# Add base Frame work / small test
# Adds a new row (row 5) of UISoundConfig data ("Click") to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 5 data --------------------------------------------------
# B5 : Id                -> 1
# C5 : NOTE               -> "Click"  (re-styled with Microsoft YaHei font)
# D5 : Resource Name      -> "Click"  (re-styled with Microsoft YaHei font)
# E5 : Priority           -> 0
# F5 : Volume             -> 1

$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "Click"
$ws.Range("D5").Value = "Click"
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1

# Give the NOTE / Resource Name cells the "Microsoft YaHei" font, as in
# the authored workbook.
$ws.Range("C5:D5").Font.Name = "Microsoft YaHei"
$ws.Range("C5:D5").Font.Size = 11

# Match the row height Excel computed for the new, taller row.
$ws.Rows.Item(5).RowHeight = 16.5

# Update the current selection to match the author's last selection.
[void]$ws.Range("G8").Select()
